$wb = $excel.ActiveWorkbook

# --- Add the new "unique column Names" sheet after the existing "PostCode" sheet ---
$postcode = $wb.Worksheets.Item("PostCode")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "unique column Names"

$newSheet.Range("A1").Value = "uniqueColumn"
$newSheet.Range("A2").Value = "EmergencyAreaIDSub"
$newSheet.Columns.Item(1).ColumnWidth = 19.5

# Header cell shading: solid fill, theme color (Background 2 / theme index 2).
# Setting Color first (reusing the existing yellow fill slot) then ThemeColor
# avoids an extra throwaway "blank" fill being recorded.
$newSheet.Range("A1").Interior.Color = 65535
$newSheet.Range("A1").Interior.ThemeColor = 4

# --- Add the "Downloadspath" column to the "PostCode" sheet ---
$postcode.Range("AM1").Value = "Downloadspath"
$postcode.Range("AM2:AM6").Value = "C:\Users\SSUPRAJA-adm\Downloads"
$postcode.Columns.Item(39).ColumnWidth = 33.5

# --- Restore "PostCode" as the active sheet/selection ---
$postcode.Activate() | Out-Null
$postcode.Range("D3").Select() | Out-Null
